$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A122:D130").EntireRow.Insert()
$v = $ws1.Range("D122").Value2
Write-Host "D122 val after insert (should be blank or old val?):" "[$v]"
